# Update attendance summary sheet: mark Real/Duplicate(Total)/Absent/Invalid
# flags as 1 for the relevant dates/columns (final code with comment, tut 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToSet = @(
    "G3", "H3",
    "D4", "E4",
    "D5", "E5",
    "D6", "E6",
    "H7",
    "H8",
    "D9", "E9",
    "H10",
    "H11",
    "D12", "E12",
    "D13", "E13",
    "H14",
    "H15",
    "H16",
    "H17",
    "H18"
)

foreach ($cellAddr in $cellsToSet) {
    $ws.Range($cellAddr).Value = 1
}
